$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 82, pushing the old rows
# 82-86 down to become rows 85-89 (formatting, including the date style
# on column D, is carried along automatically by the insert).
$ws.Rows("82:84").Insert()

# Row 82 (new): Espárragos, Banquete, Provincia de Linares
$ws.Range("A82").Value = 9
$ws.Range("B82").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C82").Value = "Metropolitana"
$ws.Range("D82").Value = 44505
$ws.Range("E82").Value = 13
$ws.Range("F82").Value = 300000000
$ws.Range("G82").Value = "Espárragos"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Banquete"
$ws.Range("J82").Value = 250
$ws.Range("K82").Value = 1400
$ws.Range("L82").Value = 1400
$ws.Range("M82").Value = 1400
$ws.Range("N82").Value = "$/kilo"
$ws.Range("O82").Value = "Provincia de Linares"
$ws.Range("P82").Value = 1400
$ws.Range("Q82").Value = 1
$ws.Range("R82").Value = "Hortaliza"

# Row 83 (new): Espárragos, Primera, Provincia de Linares
$ws.Range("A83").Value = 9
$ws.Range("B83").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C83").Value = "Metropolitana"
$ws.Range("D83").Value = 44505
$ws.Range("E83").Value = 13
$ws.Range("F83").Value = 300000000
$ws.Range("G83").Value = "Espárragos"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 160
$ws.Range("K83").Value = 1200
$ws.Range("L83").Value = 1200
$ws.Range("M83").Value = 1200
$ws.Range("N83").Value = "$/kilo"
$ws.Range("O83").Value = "Provincia de Linares"
$ws.Range("P83").Value = 1200
$ws.Range("Q83").Value = 1
$ws.Range("R83").Value = "Hortaliza"

# Row 84 (new): Espárragos, Segunda, Provincia de Linares
$ws.Range("A84").Value = 9
$ws.Range("B84").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C84").Value = "Metropolitana"
$ws.Range("D84").Value = 44505
$ws.Range("E84").Value = 13
$ws.Range("F84").Value = 300000000
$ws.Range("G84").Value = "Espárragos"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Segunda"
$ws.Range("J84").Value = 106
$ws.Range("K84").Value = 1000
$ws.Range("L84").Value = 1000
$ws.Range("M84").Value = 1000
$ws.Range("N84").Value = "$/kilo"
$ws.Range("O84").Value = "Provincia de Linares"
$ws.Range("P84").Value = 1000
$ws.Range("Q84").Value = 1
$ws.Range("R84").Value = "Hortaliza"
